# "Generate Report for Archive"
#
# Two things changed in this localization-status report:
#   1. The "Ready for handoff" status text became "In Translation"
#      (it's the shared string used by Overview!E2/F2, zh-cn!C2, de-de!C2).
#   2. The Status columns got a bit narrower (Overview!E:F, zh-cn!C, de-de!C).

$wb = $excel.ActiveWorkbook

# 1) Update the status cells that read "Ready for handoff".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# 2) Narrow the Status columns (Excel stores widths snapped to a pixel
#    grid, so 12.5 is the character-width input that lands closest to
#    the target column width).
$wsOverview.Range("E1").EntireColumn.ColumnWidth = 12.5
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 12.5
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 12.5
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 12.5
